$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.28256698319107
$ws.Range("C2").Value = 7.939454770161528
$ws.Range("D2").Value = 13.88040491818303
$ws.Range("E2").Value = 14.4491880195096
$ws.Range("G2").Value = 40.43985250332256
$ws.Range("H2").Value = 17.36511617770097
$ws.Range("I2").Value = 27.67908473172471
$ws.Range("J2").Value = 8.654041313569321
$ws.Range("K2").Value = 10.2547189018129
$ws.Range("L2").Value = 12.069061777443
$ws.Range("O2").Value = 27.9541423893907

$ws.Range("B3").Value = 13.04316414706904
$ws.Range("C3").Value = 7.910122115429834
$ws.Range("D3").Value = 13.87501188872806
$ws.Range("E3").Value = 14.47265642006962
$ws.Range("G3").Value = 40.58080681898002
$ws.Range("H3").Value = 17.41926023040808
$ws.Range("I3").Value = 27.78383042537799
$ws.Range("J3").Value = 8.665545487992295
$ws.Range("K3").Value = 10.07909805346756
$ws.Range("L3").Value = 12.06822039055388
$ws.Range("O3").Value = 28.0499834655675

$ws.Range("B4").Value = 12.89567385034402
$ws.Range("C4").Value = 7.892038584171324
$ws.Range("D4").Value = 13.87416933915733
$ws.Range("E4").Value = 14.48901388403092
$ws.Range("G4").Value = 40.676891712589
$ws.Range("H4").Value = 17.45487594781716
$ws.Range("I4").Value = 27.85257317439932
$ws.Range("J4").Value = 8.67300260886172
$ws.Range("K4").Value = 9.970830624051793
$ws.Range("L4").Value = 12.06925179273881
$ws.Range("O4").Value = 28.11367870140365

$ws.Range("B5").Value = 12.83552135808816
$ws.Range("C5").Value = 7.884653308974393
$ws.Range("D5").Value = 13.87444871014071
$ws.Range("E5").Value = 14.49616988711887
$ws.Range("G5").Value = 40.71843982846661
$ws.Range("H5").Value = 17.46998633315476
$ws.Range("I5").Value = 27.8817004299226
$ws.Range("J5").Value = 8.676140683826816
$ws.Range("K5").Value = 9.926653758190835
$ws.Range("L5").Value = 12.07006266500565
$ws.Range("O5").Value = 28.14085366158757

$ws.Range("B6").Value = 12.82553236364953
$ws.Range("C6").Value = 7.883426084855828
$ws.Range("D6").Value = 13.87453276013251
$ws.Range("E6").Value = 14.49738775451335
$ws.Range("G6").Value = 40.72548320438981
$ws.Range("H6").Value = 17.47253145382525
$ws.Range("I6").Value = 27.88660428682187
$ws.Range("J6").Value = 8.676667760916967
$ws.Range("K6").Value = 9.919316332025932
$ws.Range("L6").Value = 12.07022093491788
$ws.Range("O6").Value = 28.14543962961628

$ws.Range("B7").Value = 12.89486271010005
$ws.Range("C7").Value = 7.891939045881108
$ws.Range("D7").Value = 13.87417058317963
$ws.Range("E7").Value = 14.4891084070457
$ws.Range("G7").Value = 40.67744236468079
$ws.Range("H7").Value = 17.45507731515005
$ws.Range("I7").Value = 27.85296148348509
$ws.Range("J7").Value = 8.673044527847102
$ws.Range("K7").Value = 9.970235000183303
$ws.Range("L7").Value = 12.0692611454109
$ws.Range("O7").Value = 28.11404025894228

$ws.Range("B8").Value = 13.20016877713917
$ws.Range("C8").Value = 7.929357443311456
$ws.Range("D8").Value = 13.87803424249683
$ws.Range("E8").Value = 14.45687592903939
$ws.Range("G8").Value = 40.48647046729285
$ws.Range("H8").Value = 17.38329321939847
$ws.Range("I8").Value = 27.71428218900714
$ws.Range("J8").Value = 8.657926452038842
$ws.Range("K8").Value = 10.19428641922056
$ws.Range("L8").Value = 12.06845119175082
$ws.Range("O8").Value = 27.98618178099806

$ws.Range("B9").Value = 13.79169869474578
$ws.Range("C9").Value = 8.002053698064858
$ws.Range("D9").Value = 13.90510371824857
$ws.Range("E9").Value = 14.40910333820491
$ws.Range("G9").Value = 40.18790663368415
$ws.Range("H9").Value = 17.2613167835029
$ws.Range("I9").Value = 27.47744981517834
$ws.Range("J9").Value = 8.631389262297324
$ws.Range("K9").Value = 10.62794470391239
$ws.Range("L9").Value = 12.07908479349429
$ws.Range("O9").Value = 27.77394554412115

$ws.Range("B10").Value = 14.21774272729431
$ws.Range("C10").Value = 8.054914809681515
$ws.Range("D10").Value = 13.93672933227095
$ws.Range("E10").Value = 14.38338676928351
$ws.Range("G10").Value = 40.0151683522322
$ws.Range("H10").Value = 17.18312738332078
$ws.Range("I10").Value = 27.32482873929535
$ws.Range("J10").Value = 8.613769718918425
$ws.Range("K10").Value = 10.94016381175291
$ws.Range("L10").Value = 12.09426115010067
$ws.Range("O10").Value = 27.64151418788149

$ws.Range("B11").Value = 14.40887384808009
$ws.Range("C11").Value = 8.078814824334758
$ws.Range("D11").Value = 13.95362786745478
$ws.Range("E11").Value = 14.37371796755204
$ws.Range("G11").Value = 39.9467721325549
$ws.Range("H11").Value = 17.1500321458583
$ws.Range("I11").Value = 27.26003433946283
$ws.Range("J11").Value = 8.606157915565559
$ws.Range("K11").Value = 11.08023199165505
$ws.Range("L11").Value = 12.10274328413097
$ws.Range("O11").Value = 27.58637697122451

$ws.Range("B12").Value = 14.48080016194234
$ws.Range("C12").Value = 8.087841525053225
$ws.Range("D12").Value = 13.96038438144395
$ws.Range("E12").Value = 14.37034785032212
$ws.Range("G12").Value = 39.92234146001638
$ws.Range("H12").Value = 17.13785511640294
$ws.Range("I12").Value = 27.23616444530266
$ws.Range("J12").Value = 8.603333242114452
$ws.Range("K12").Value = 11.13294440824395
$ws.Range("L12").Value = 12.10618019457312
$ws.Range("O12").Value = 27.56623286811017

$ws.Range("B13").Value = 14.46533060579939
$ws.Range("C13").Value = 8.085898559993757
$ws.Range("D13").Value = 13.95891341345631
$ws.Range("E13").Value = 14.37106072274885
$ws.Range("G13").Value = 39.92753760810873
$ws.Range("H13").Value = 17.14046185441896
$ws.Range("I13").Value = 27.24127561786953
$ws.Range("J13").Value = 8.603939021684285
$ws.Range("K13").Value = 11.12160716162079
$ws.Range("L13").Value = 12.10543002550179
$ws.Range("O13").Value = 27.57053855100851

$ws.Range("B14").Value = 14.41480066969183
$ws.Range("C14").Value = 8.079557946728295
$ws.Range("D14").Value = 13.9541765898855
$ws.Range("E14").Value = 14.37343487277156
$ws.Range("G14").Value = 39.94473272735056
$ws.Range("H14").Value = 17.1490232127058
$ws.Range("I14").Value = 27.25805719426303
$ws.Range("J14").Value = 8.605924372036299
$ws.Range("K14").Value = 11.08457551134453
$ws.Range("L14").Value = 12.10302154549068
$ws.Range("O14").Value = 27.58470496508846

$ws.Range("B15").Value = 14.3837890147129
$ws.Range("C15").Value = 8.075670971838861
$ws.Range("D15").Value = 13.95132157520083
$ws.Range("E15").Value = 14.37492701750543
$ws.Range("G15").Value = 39.9554567525167
$ws.Range("H15").Value = 17.15431357003245
$ws.Range("I15").Value = 27.26842317243842
$ws.Range("J15").Value = 8.607147969716905
$ws.Range("K15").Value = 11.0618484610766
$ws.Range("L15").Value = 12.10157550882554
$ws.Range("O15").Value = 27.59347806806434

$ws.Range("B16").Value = 14.20519226212214
$ws.Range("C16").Value = 8.053349711177466
$ws.Range("D16").Value = 13.93567520120025
$ws.Range("E16").Value = 14.38405944002812
$ws.Range("G16").Value = 40.01984356751102
$ws.Range("H16").Value = 17.18533998507696
$ws.Range("I16").Value = 27.32915644842804
$ws.Range("J16").Value = 8.614275263693587
$ws.Range("K16").Value = 10.93096654239651
$ws.Range("L16").Value = 12.09373841030902
$ws.Range("O16").Value = 27.64522036380598

$ws.Range("B17").Value = 14.09489504132886
$ws.Range("C17").Value = 8.039616895126809
$ws.Range("D17").Value = 13.92671740882872
$ws.Range("E17").Value = 14.39018133209605
$ws.Range("G17").Value = 40.06195474603233
$ws.Range("H17").Value = 17.20500701704411
$ws.Range("I17").Value = 27.36760115029138
$ws.Range("J17").Value = 8.61875076676545
$ws.Range("K17").Value = 10.85013860365575
$ws.Range("L17").Value = 12.08933343395327
$ws.Range("O17").Value = 27.67827109320505

$ws.Range("B18").Value = 14.03120660294727
$ws.Range("C18").Value = 8.031704484770467
$ws.Range("D18").Value = 13.92180180187692
$ws.Range("E18").Value = 14.39389357796045
$ws.Range("G18").Value = 40.08713419454884
$ws.Range("H18").Value = 17.21655180787624
$ws.Range("I18").Value = 27.39014966799525
$ws.Range("J18").Value = 8.621362944478291
$ws.Range("K18").Value = 10.8034664744219
$ws.Range("L18").Value = 12.0869484991947
$ws.Range("O18").Value = 27.69776159340999

$ws.Range("B19").Value = 14.00960219061555
$ws.Range("C19").Value = 8.029023212220556
$ws.Range("D19").Value = 13.92017822211798
$ws.Range("E19").Value = 14.39518331980616
$ws.Range("G19").Value = 40.09582393481118
$ws.Range("H19").Value = 17.22050067356072
$ws.Range("I19").Value = 27.39785911638367
$ws.Range("J19").Value = 8.622253915164228
$ws.Range("K19").Value = 10.78763427570941
$ws.Range("L19").Value = 12.08616659953189
$ws.Range("O19").Value = 27.7044432580554

$ws.Range("B20").Value = 14.1066625640199
$ws.Range("C20").Value = 8.041080209250199
$ws.Range("D20").Value = 13.92764651232848
$ws.Range("E20").Value = 14.38950987298091
$ws.Range("G20").Value = 40.05737272943546
$ws.Range("H20").Value = 17.20288933242441
$ws.Range("I20").Value = 27.36346350703174
$ws.Range("J20").Value = 8.618270412226972
$ws.Range("K20").Value = 10.85876205826234
$ws.Range("L20").Value = 12.08978697669293
$ws.Range("O20").Value = 27.67470304058708

$ws.Range("B21").Value = 14.42965528018552
$ws.Range("C21").Value = 8.081421003043758
$ws.Range("D21").Value = 13.9555582417715
$ws.Range("E21").Value = 14.37272962831208
$ws.Range("G21").Value = 39.93964218310419
$ws.Range("H21").Value = 17.14649889038547
$ws.Range("I21").Value = 27.25310995393628
$ws.Range("J21").Value = 8.605339661214755
$ws.Range("K21").Value = 11.09546187331588
$ws.Range("L21").Value = 12.10372288743728
$ws.Range("O21").Value = 27.58052398696243

$ws.Range("B22").Value = 14.63809358040222
$ws.Range("C22").Value = 8.107646591240774
$ws.Range("D22").Value = 13.97588141723535
$ws.Range("E22").Value = 14.36346004964262
$ws.Range("G22").Value = 39.87126593062469
$ws.Range("H22").Value = 17.11171601052914
$ws.Range("I22").Value = 27.18487149065538
$ws.Range("J22").Value = 8.597225167874081
$ws.Range("K22").Value = 11.24822515854089
$ws.Range("L22").Value = 12.11414071857645
$ws.Range("O22").Value = 27.52325801700755

$ws.Range("B23").Value = 14.52710996462469
$ws.Range("C23").Value = 8.093663125877388
$ws.Range("D23").Value = 13.96484545908118
$ws.Range("E23").Value = 14.3682523215307
$ws.Range("G23").Value = 39.90697408794302
$ws.Range("H23").Value = 17.13009082923935
$ws.Range("I23").Value = 27.22093624568007
$ws.Range("J23").Value = 8.601525321007303
$ws.Range("K23").Value = 11.16688412222233
$ws.Range("L23").Value = 12.10846138027381
$ws.Range("O23").Value = 27.55342954116279

$ws.Range("B24").Value = 14.10134332507075
$ws.Range("C24").Value = 8.040418698415374
$ws.Range("D24").Value = 13.92722573439629
$ws.Range("E24").Value = 14.38981283940086
$ws.Range("G24").Value = 40.059441240625
$ws.Range("H24").Value = 17.20384599640828
$ws.Range("I24").Value = 27.36533274582524
$ws.Range("J24").Value = 8.618487458541804
$ws.Range("K24").Value = 10.85486402341288
$ws.Range("L24").Value = 12.08958147028397
$ws.Range("O24").Value = 27.67631463342031

$ws.Range("B25").Value = 13.63288399592555
$ws.Range("C25").Value = 7.982474060333971
$ws.Range("D25").Value = 13.8957076326131
$ws.Range("E25").Value = 14.42037730384554
$ws.Range("G25").Value = 40.26051226097213
$ws.Range("H25").Value = 17.2923059103961
$ws.Range("I25").Value = 27.53776279174367
$ws.Range("J25").Value = 8.638237290889466
$ws.Range("K25").Value = 10.62794470391239
$ws.Range("L25").Value = 12.07490815872328
$ws.Range("O25").Value = 27.82723672592185

